$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update header labels and municipality/state text cleanup (title-casing "de/del/la/las/los/el" particles, plus a couple of one-off fixes)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B5').Value = 'Rincón De Romos'
$ws.Range('B19').Value = 'Amatenango De La Frontera'
$ws.Range('B21').Value = 'Benemérito De Las Américas'
$ws.Range('B26').Value = 'Comitán De Domínguez'
$ws.Range('B48').Value = 'San Cristóbal De Las Casas'
$ws.Range('B73').Value = 'Hidalgo Del Parral'
$ws.Range('B76').Value = 'San Francisco De Borja'
$ws.Range('A90').Value = 'Ciudad De México'
$ws.Range('B107').Value = 'San Juan De Guadalupe'
$ws.Range('A113').Value = 'Estado De México'
$ws.Range('B113').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B114').Value = 'Almoloya De Alquisiras'
$ws.Range('B117').Value = 'Atizapán De Zaragoza'
$ws.Range('B122').Value = 'Ecatepec De Morelos'
$ws.Range('B125').Value = 'Ixtapan De La Sal'
$ws.Range('B131').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B139').Value = 'Tlalnepantla De Baz'
$ws.Range('B142').Value = 'Villa De Allende'
$ws.Range('B143').Value = 'Villa Del Carbón'
$ws.Range('A148').Value = 'Guanajuato'
$ws.Range('B150').Value = 'Apaseo El Alto'
$ws.Range('B151').Value = 'Apaseo El Grande'
$ws.Range('B156').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B159').Value = 'Jaral Del Progreso'
$ws.Range('B165').Value = 'San Diego De La Unión'
$ws.Range('B167').Value = 'San Francisco Del Rincón'
$ws.Range('B169').Value = 'San Luis De La Paz'
$ws.Range('B173').Value = 'Valle De Santiago'
$ws.Range('B178').Value = 'Acapulco De Juárez'
$ws.Range('B179').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B180').Value = 'Alcozauca De Guerero'
$ws.Range('B182').Value = 'Atoyac De Álvarez'
$ws.Range('B183').Value = 'Ayutla De Los Libres'
$ws.Range('B185').Value = 'Chilapa De Álvarez'
$ws.Range('B186').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B187').Value = 'Coyuca De Benítez'
$ws.Range('B189').Value = 'Cutzamala De Pinzón'
$ws.Range('B191').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B193').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B201').Value = 'Técpan De Galeana'
$ws.Range('B202').Value = 'Tixtla De Guerero'
$ws.Range('B213').Value = 'Cuautepec De Hinojosa'
$ws.Range('B219').Value = 'Pachuca De Soto'
$ws.Range('B222').Value = 'Tenango De Doria'
$ws.Range('B223').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B226').Value = 'Tula De Allende'
$ws.Range('B228').Value = 'Zacualtipán De Ángeles'
$ws.Range('B234').Value = 'Autlán De Navarro'
$ws.Range('B243').Value = 'Lagos De Moreno'
$ws.Range('B249').Value = 'San Miguel El Alto'
$ws.Range('B250').Value = 'Tamazula De Gordiano'
$ws.Range('B252').Value = 'Tepatitlán De Morelos'
$ws.Range('B294').Value = 'Coatlán Del Río'
$ws.Range('B300').Value = 'Tlaltizapán De Zapata'
$ws.Range('B314').Value = 'Montemorelos'
$ws.Range('B316').Value = 'San Nicolás De Los Garza'
$ws.Range('B319').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B323').Value = 'Coicoyán De Las Flores'
$ws.Range('B325').Value = 'Ixtlán De Juárez'
$ws.Range('B327').Value = 'Oaxaca De Juárez'
$ws.Range('B328').Value = 'Putla Villa De Guerero'
$ws.Range('B331').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B366').Value = 'Tanetze De Zaragoza'
$ws.Range('B367').Value = 'Tataltepec De Valdés'
$ws.Range('B368').Value = 'Teotitlán De Flores Magón'
$ws.Range('B369').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B381').Value = 'Palmar De Bravo'
$ws.Range('B386').Value = 'Tecali De Herrera'
$ws.Range('B388').Value = 'Tepexi De Rodríguez'
$ws.Range('B389').Value = 'Tetela De Ocampo'
$ws.Range('B397').Value = 'Xayacatlán De Bravo'
$ws.Range('B406').Value = 'Cadereyta De Montes'
$ws.Range('B409').Value = 'Jalpan De Serra'
$ws.Range('B410').Value = 'Pinal De Amoles'
$ws.Range('B413').Value = 'San Juan Del Río'
$ws.Range('B420').Value = 'Ciudad Del Maíz'
$ws.Range('B426').Value = 'Mexquitic De Carmona'
$ws.Range('B430').Value = 'San Ciro De Acosta'
$ws.Range('B434').Value = 'Santa María Del Río'
$ws.Range('B440').Value = 'Villa De Arriaga'
$ws.Range('B441').Value = 'Villa De Ramos'
$ws.Range('B442').Value = 'Villa De Reyes'
$ws.Range('B464').Value = 'Jalpa De Méndez'
$ws.Range('B487').Value = 'Soto La Marina'
$ws.Range('B498').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B506').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B508').Value = 'Amatlán De Los Reyes'
$ws.Range('B514').Value = 'Boca Del Río'
$ws.Range('B519').Value = 'Cosamaloapan De Carpio'
$ws.Range('B528').Value = 'Hueyapan De Ocampo'
$ws.Range('B529').Value = 'Ignacio De La Llave'
$ws.Range('B532').Value = 'Ixhuatlán Del Café'
$ws.Range('B533').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B540').Value = 'Lerdo De Tejada'
$ws.Range('B542').Value = 'Martínez De La Torre'
$ws.Range('B546').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B558').Value = 'Sayula De Alemán'
$ws.Range('B559').Value = 'Soledad De Doblado'
$ws.Range('B561').Value = 'Tatahuicapan De Juárez'
$ws.Range('B577').Value = 'Vega De Alatorre'
$ws.Range('B588').Value = 'Cañitas De Felipe Pescador'

# 2) Remove trailing metadata/footer rows 604-608 (sample size, source, author, date notes)
$ws.Range("A604:D608").EntireRow.Delete() | Out-Null
